# SVM, LR, DNN comparison
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "y"

# --- Existing rows, now paired with counts ---
$ws.Range("A2").Value = "How to cold start a meter using 1132?"
$ws.Range("B2").Style = "Normal"
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = "What is RadioShop?"
$ws.Range("B3").Value = 2

# --- New rows ---
$ws.Range("A4").Value = "How can I install Command Center?"
$ws.Range("B4").Value = 12

$ws.Range("A5").Value = "Send Seed files"
$ws.Range("B5").Value = 17

$ws.Range("A6").Value = "How can I perform a PCA test on a meter?"
$ws.Range("B6").Value = 25

$ws.Range("A7").Value = "How to install RadioShop"
$ws.Range("B7").Value = 2

$ws.Range("A8").Value = "where do I find DCW Folder?"
$ws.Range("B8").Value = 15

# --- Column width (A: ~35.37 chars wide) ---
$ws.Columns.Item(1).ColumnWidth = 34.534

# --- Selection ---
[void]$ws.Range("G21").Select()
